$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Fix the title text: "92 genes" -> "79 genes" (A1), preserving
#    the italic "Lactobacillus plantarum" run at the end.
# ---------------------------------------------------------------
$titleCell = $ws.Range("A1")
$titleCell.Characters(49, 2).Text = "79"
$italicRun = $titleCell.Characters(118, 24)
$italicRun.Font.Italic = $true

# ---------------------------------------------------------------
# 2. Update the data table (rows 4-17): new order, frequencies and
#    percentages based on 79 genes / 152 total occurrences instead
#    of the old 175 total occurrences.
# ---------------------------------------------------------------
$ws.Range("A4").Value = "Amino acid metabolism"
$ws.Range("B4").Value = 42
$ws.Range("C4").Value = 27.631578947368421

$ws.Range("A5").Value = "Carbohydrate metabolism"
$ws.Range("B5").Value = 40
$ws.Range("C5").Value = 26.315789473684209

$ws.Range("A6").Value = "Nucleotide metabolism"
$ws.Range("B6").Value = 19
$ws.Range("C6").Value = 12.5

$ws.Range("A7").Value = "Metabolism of cofactors and vitamins"
$ws.Range("B7").Value = 11
$ws.Range("C7").Value = 7.2368421052631575

$ws.Range("A8").Value = "Energy metabolism"
$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 5.2631578947368425

$ws.Range("A9").Value = "Metabolism of other amino acids"
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 5.2631578947368425

$ws.Range("A10").Value = "Biosynthesis of other secondary metabolites"
$ws.Range("B10").Value = 7
$ws.Range("C10").Value = 4.6052631578947372

$ws.Range("A11").Value = "Glycan biosynthesis and metabolism"
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 3.2894736842105261

$ws.Range("A12").Value = "Lipid metabolism"
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 2.6315789473684212

$ws.Range("A13").Value = "Metabolism of terpenoids and polyketides"
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 2.6315789473684212

$ws.Range("A14").Value = "Quorum sensing"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 0.65789473684210531

$ws.Range("A15").Value = "RNA degradation"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 0.65789473684210531

$ws.Range("A16").Value = "Two-component system"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 0.65789473684210531

$ws.Range("A17").Value = "Xenobiotics biodegradation and metabolism"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 0.65789473684210531

# ---------------------------------------------------------------
# 3. Strip the (now unwanted) centered formatting from columns A
#    and B, and from column C -- column C keeps its "0.00" number
#    format but loses the center alignment.
# ---------------------------------------------------------------
$ws.Range("A4:B17").Style = "Normal"
$ws.Range("C4:C17").Style = "Normal"
$ws.Range("C4:C17").NumberFormat = "0.00"

# ---------------------------------------------------------------
# 4. Widen column C slightly (12.15625 -> ~13.15625 chars).
# ---------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 12.26
